$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 14:05"

# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 1528179
$ws.Cells.Item(4, 3).Value = 515
$ws.Cells.Item(4, 4).Value = 346389
$ws.Cells.Item(4, 5).Value = 1090802
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 10
$ws.Cells.Item(4, 8).Value = 90988

# Row 13: Iran
$ws.Cells.Item(13, 1).Value = "Iran"
$ws.Cells.Item(13, 2).Value = 122492
$ws.Cells.Item(13, 3).Value = 2294
$ws.Cells.Item(13, 4).Value = 95661
$ws.Cells.Item(13, 5).Value = 19774
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 69
$ws.Cells.Item(13, 8).Value = 7057

# Row 21: Paises Bajos
$ws.Cells.Item(21, 1).Value = "Paises Bajos"
$ws.Cells.Item(21, 2).Value = 44141
$ws.Cells.Item(21, 3).Value = 146
$ws.Cells.Item(21, 4).Value = 0
$ws.Cells.Item(21, 5).Value = 0
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 14
$ws.Cells.Item(21, 8).Value = 5694

# Row 28: Suecia
$ws.Cells.Item(28, 1).Value = "Suecia"
$ws.Cells.Item(28, 2).Value = 30377
$ws.Cells.Item(28, 3).Value = 234
$ws.Cells.Item(28, 4).Value = 4971
$ws.Cells.Item(28, 5).Value = 21708
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).Value = 19
$ws.Cells.Item(28, 8).Value = 3698

# Row 29: Portugal
$ws.Cells.Item(29, 1).Value = "Portugal"
$ws.Cells.Item(29, 2).Value = 29209
$ws.Cells.Item(29, 3).Value = 173
$ws.Cells.Item(29, 4).Value = 6430
$ws.Cells.Item(29, 5).Value = 21548
$ws.Cells.Item(29, 6).Value = 0
$ws.Cells.Item(29, 7).Value = 13
$ws.Cells.Item(29, 8).Value = 1231

# Row 41: Kuwait
$ws.Cells.Item(41, 1).Value = "Kuwait"
$ws.Cells.Item(41, 2).Value = 15691
$ws.Cells.Item(41, 3).Value = 841
$ws.Cells.Item(41, 4).Value = 4339
$ws.Cells.Item(41, 5).Value = 11234
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(41, 7).Value = 6
$ws.Cells.Item(41, 8).Value = 118

# Row 42: Colombia
$ws.Cells.Item(42, 1).Value = "Colombia"
$ws.Cells.Item(42, 2).Value = 15574
$ws.Cells.Item(42, 3).Value = 0
$ws.Cells.Item(42, 4).Value = 3751
$ws.Cells.Item(42, 5).Value = 11249
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(42, 7).Value = 0
$ws.Cells.Item(42, 8).Value = 574

# Row 43: Sudafrica
$ws.Cells.Item(43, 1).Value = "Sudafrica"
$ws.Cells.Item(43, 2).Value = 15515
$ws.Cells.Item(43, 3).Value = 0
$ws.Cells.Item(43, 4).Value = 7006
$ws.Cells.Item(43, 5).Value = 8245
$ws.Cells.Item(43, 6).Value = 0
$ws.Cells.Item(43, 7).Value = 0
$ws.Cells.Item(43, 8).Value = 264

# Row 55: Australia
$ws.Cells.Item(55, 1).Value = "Australia"
$ws.Cells.Item(55, 2).Value = 7060
$ws.Cells.Item(55, 3).Value = 15
$ws.Cells.Item(55, 4).Value = 6392
$ws.Cells.Item(55, 5).Value = 569
$ws.Cells.Item(55, 6).Value = 0
$ws.Cells.Item(55, 7).Value = 1
$ws.Cells.Item(55, 8).Value = 99

# Row 75: Uzbekistan
$ws.Cells.Item(75, 1).Value = "Uzbekistan"
$ws.Cells.Item(75, 2).Value = 2779
$ws.Cells.Item(75, 3).Value = 26
$ws.Cells.Item(75, 4).Value = 2293
$ws.Cells.Item(75, 5).Value = 473
$ws.Cells.Item(75, 6).Value = 0
$ws.Cells.Item(75, 7).Value = 1
$ws.Cells.Item(75, 8).Value = 13

# Row 112: Mali
$ws.Cells.Item(112, 1).Value = "Mali"
$ws.Cells.Item(112, 2).Value = 874
$ws.Cells.Item(112, 3).Value = 14
$ws.Cells.Item(112, 4).Value = 512
$ws.Cells.Item(112, 5).Value = 310
$ws.Cells.Item(112, 6).Value = 0
$ws.Cells.Item(112, 7).Value = 0
$ws.Cells.Item(112, 8).Value = 52

# Row 113: Costa Rica
$ws.Cells.Item(113, 1).Value = "Costa Rica"
$ws.Cells.Item(113, 2).Value = 863
$ws.Cells.Item(113, 3).Value = 0
$ws.Cells.Item(113, 4).Value = 565
$ws.Cells.Item(113, 5).Value = 288
$ws.Cells.Item(113, 6).Value = 0
$ws.Cells.Item(113, 7).Value = 0
$ws.Cells.Item(113, 8).Value = 10

# Row 116: Zambia
$ws.Cells.Item(116, 1).Value = "Zambia"
$ws.Cells.Item(116, 2).Value = 761
$ws.Cells.Item(116, 3).Value = 8
$ws.Cells.Item(116, 4).Value = 192
$ws.Cells.Item(116, 5).Value = 562
$ws.Cells.Item(116, 6).Value = 0
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 7

# Row 117: Principado de Andorra
$ws.Cells.Item(117, 1).Value = "Principado de Andorra"
$ws.Cells.Item(117, 2).Value = 761
$ws.Cells.Item(117, 3).Value = 0
$ws.Cells.Item(117, 4).Value = 617
$ws.Cells.Item(117, 5).Value = 93
$ws.Cells.Item(117, 6).Value = 0
$ws.Cells.Item(117, 7).Value = 0
$ws.Cells.Item(117, 8).Value = 51
